# Update the "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
# Both sheets carry the same event list in rows 2-32; F8 and F11 have
# slightly different source values between the two sheets that converge
# to different (F8) or the same (F11) totals after the refresh.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# row -> new value, shared between the two sheets
$common = @{
    5  = 39
    6  = 513
    7  = 44
    10 = 85
    13 = 268
    15 = 81
    16 = 14
    18 = 2826
    19 = 43
    20 = 394
    21 = 15
    23 = 59
    25 = 51
    26 = 5
    28 = 40
    29 = 184
    30 = 254
    31 = 1612
    32 = 227
}

# rows where the two sheets diverge
$perSheet = @{
    "展览"   = @{ 8 = 1954; 11 = 4112 }
    "全部类型" = @{ 8 = 1955; 11 = 4112 }
}

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    foreach ($row in $common.Keys) {
        $ws.Cells.Item($row, 6).Value = $common[$row]
    }

    $special = $perSheet[$name]
    foreach ($row in $special.Keys) {
        $ws.Cells.Item($row, 6).Value = $special[$row]
    }
}
